$wb = $excel.ActiveWorkbook

# --- Add the new "data to plot" sheet as the first sheet ---
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "data to plot"
$newSheet.Move($wb.Worksheets.Item(1))
$plotSheet = $wb.Worksheets.Item("data to plot")

$zinc = $wb.Worksheets.Item("Zinc")
$germanium = $wb.Worksheets.Item("Germanium")
$co2 = $wb.Worksheets.Item("CO2")

# --- Add "Substance" / substance-name header row (row 5) to each data sheet ---
# Order chosen to reproduce the shared-string table ordering of the target file.
$zinc.Range("D5").Value = "Substance"
$plotSheet.Range("B1").Value = "Multiple_of_sustainable"
$germanium.Range("E5").Value = "Germanium"
$zinc.Range("E5").Value = "Zinc"
$co2.Range("E5").Value = "CO2"
$plotSheet.Range("A1").Value = "Resource"

$germanium.Range("D5").Value = "Substance"
$co2.Range("D5").Value = "Substance"

# CO2's new header cells carry an explicit black font color (new style in target file)
$co2.Range("D5:E5").Font.Color = 0

# --- Populate the "data to plot" sheet with formulas referencing each data sheet ---
$plotSheet.Range("A2").Formula = "=Zinc!E5"
$plotSheet.Range("B2").Formula = "=Zinc!E24"

$plotSheet.Range("A3").Formula = "=Germanium!E5"
$plotSheet.Range("B3").Formula = "=Germanium!E24"

$plotSheet.Range("A4").Formula = "='CO2'!E5"
$plotSheet.Range("B4").Formula = "='CO2'!E23"

# --- Update selections to match where the authors were working ---
$zinc.Activate()
$zinc.Range("D5:E5").Select()

$germanium.Activate()
$germanium.Range("E5").Select()

$co2.Activate()
$co2.Range("E6").Select()

$plotSheet.Activate()
$plotSheet.Range("A1").Select()

Write-Host "done"
